# Fruta / hortaliza, semanal
# Updates the weekly price records (rows 2-19) on the active sheet with the
# refreshed values (the underlying records were reshuffled/updated across
# the existing date rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2;  D=44410; H="Sin especificar"; I="Primera"; J=100; K=5500; L=6000; M=5750; P=5750 },
    @{ Row=3;  D=44636; H="Americana (o)";   I="Primera"; J=60;  K=8000; L=9000; M=8500; P=8500 },
    @{ Row=4;  D=44259; H="Sin especificar"; I="Primera"; J=80;  K=4000; L=4500; M=4250; P=4250 },
    @{ Row=5;  D=44539; H="Americana (o)";   I="Primera"; J=160; K=6500; L=7000; M=6750; P=6750 },
    @{ Row=6;  D=44699; H="Sin especificar"; I="Primera"; J=50;  K=9000; L=9500; M=9250; P=9250 },
    @{ Row=7;  D=44263; H="Sin especificar"; I="Primera"; J=100; K=7000; L=8000; M=7500; P=7500 },
    @{ Row=8;  D=44309; H="Sin especificar"; I="Primera"; J=50;  K=8000; L=9000; M=8500; P=8500 },
    @{ Row=9;  D=44945; H="Sin especificar"; I="Primera"; J=45;  K=6000; L=7000; M=6444; P=6444 },
    @{ Row=10; D=44804; H="Sin especificar"; I="Primera"; J=60;  K=5500; L=6000; M=5750; P=5750 },
    @{ Row=11; D=44789; H="Sin especificar"; I="Primera"; J=80;  K=5000; L=6000; M=5500; P=5500 },
    @{ Row=12; D=44497; H="Sin especificar"; I="Primera"; J=160; K=5000; L=6000; M=5500; P=5500 },
    @{ Row=13; D=44764; H="Americana (o)";   I="Primera"; J=100; K=7000; L=8000; M=7500; P=7500 },
    @{ Row=14; D=44414; H="Sin especificar"; I="Primera"; J=100; K=6000; L=7000; M=6500; P=6500 },
    @{ Row=15; D=44281; H="Sin especificar"; I="Primera"; J=100; K=5000; L=6000; M=5500; P=5500 },
    @{ Row=16; D=44559; H="Americana (o)";   I="Primera"; J=100; K=5000; L=6000; M=5500; P=5500 },
    @{ Row=17; D=44371; H="Sin especificar"; I="Primera"; J=80;  K=7000; L=8000; M=7375; P=7375 },
    @{ Row=18; D=44575; H="Sin especificar"; I="Primera"; J=160; K=6500; L=7000; M=6750; P=6750 },
    @{ Row=19; D=44253; H="Americana (o)";   I="Segunda"; J=100; K=4000; L=4500; M=4250; P=4250 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("D$n").Value = $r.D
    $ws.Range("H$n").Value = $r.H
    $ws.Range("I$n").Value = $r.I
    $ws.Range("J$n").Value = $r.J
    $ws.Range("K$n").Value = $r.K
    $ws.Range("L$n").Value = $r.L
    $ws.Range("M$n").Value = $r.M
    $ws.Range("P$n").Value = $r.P
}
